# Hourlog.xlsx: append a new tutoring-session row (2017-11-30, 2 hours) to
# the "DI" sheet, matching the formatting of the row above, and move the
# active-cell selection onto the new date cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row 13 — this pushes formatting down from row 12 so the
# new A13/B13 cells pick up the same date / centered-number styles used
# by the rest of the log, without minting any new number formats.
$ws.Rows(13).Insert() | Out-Null

# New log entry: date serial 43069 = 2017-11-30, 2 hours logged.
$ws.Range("A13").Value = 43069
$ws.Range("B13").Value = 2

# Formula in E2 (E1-SUM(B:B)) recalculates automatically once B13 is set.

# Match the saved selection state: active cell on the new date cell.
$ws.Range("A13").Select() | Out-Null
